$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 127.134477
$ws.Range("H2").Value = 381.4034310000001
$ws.Range("I2").Value = 0.7397224744116364
$ws.Range("J2").Value = 0.7397224744116365
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.03403033333333
$ws.Range("N2").Value = 63.102091
$ws.Range("O2").Value = 0.1030737782409588
$ws.Range("P2").Value = 0.1030737782409588
$ws.Range("Q2").Value = 2674.150445630469
$ws.Range("R2").Value = 24067.35401067423
$ws.Range("S2").Value = 0.07624599028735835
$ws.Range("T2").Value = 0.07624599028735836

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 127.134477
$ws.Range("H3").Value = 381.4034310000001
$ws.Range("I3").Value = 0.7397224744116364
$ws.Range("J3").Value = 0.7397224744116365
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.853217
$ws.Range("N3").Value = 5.559651000000001
$ws.Range("O3").Value = 0.009081382648177616
$ws.Range("P3").Value = 0.009081382648177616
$ws.Range("Q3").Value = 235.607774062509
$ws.Range("R3").Value = 2120.469966562582
$ws.Range("S3").Value = 0.006717702843588845
$ws.Range("T3").Value = 0.006717702843588847

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 127.134477
$ws.Range("H4").Value = 381.4034310000001
$ws.Range("I4").Value = 0.7397224744116364
$ws.Range("J4").Value = 0.7397224744116365
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 179.590866
$ws.Range("N4").Value = 538.772598
$ws.Range("O4").Value = 0.8800552629635879
$ws.Range("P4").Value = 0.8800552629635878
$ws.Range("Q4").Value = 22832.19082288709
$ws.Range("R4").Value = 205489.7174059838
$ws.Range("S4").Value = 0.6509966567384087
$ws.Range("T4").Value = 0.6509966567384087

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 127.134477
$ws.Range("H5").Value = 381.4034310000001
$ws.Range("I5").Value = 0.7397224744116364
$ws.Range("J5").Value = 0.7397224744116365
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.589601
$ws.Range("N5").Value = 4.768803
$ws.Range("O5").Value = 0.007789576147275677
$ws.Range("P5").Value = 0.007789576147275675
$ws.Range("Q5").Value = 202.093091773677
$ws.Range("R5").Value = 1818.837825963093
$ws.Range("S5").Value = 0.005762124542280625
$ws.Range("T5").Value = 0.005762124542280625

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.89437833333333
$ws.Range("H6").Value = 89.68313499999999
$ws.Range("I6").Value = 0.1739382112039596
$ws.Range("J6").Value = 0.1739382112039596
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 21.03403033333333
$ws.Range("N6").Value = 63.102091
$ws.Range("O6").Value = 0.1030737782409588
$ws.Range("P6").Value = 0.1030737782409588
$ws.Range("Q6").Value = 628.7992606594761
$ws.Range("R6").Value = 5659.193345935285
$ws.Range("S6").Value = 0.01792846860926599
$ws.Range("T6").Value = 0.01792846860926599

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.89437833333333
$ws.Range("H7").Value = 89.68313499999999
$ws.Range("I7").Value = 0.1739382112039596
$ws.Range("J7").Value = 0.1739382112039596
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.853217
$ws.Range("N7").Value = 5.559651000000001
$ws.Range("O7").Value = 0.009081382648177616
$ws.Range("P7").Value = 0.009081382648177616
$ws.Range("Q7").Value = 55.400770131765
$ws.Range("R7").Value = 498.606931185885
$ws.Range("S7").Value = 0.001579599453082692
$ws.Range("T7").Value = 0.001579599453082692

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.89437833333333
$ws.Range("H8").Value = 89.68313499999999
$ws.Range("I8").Value = 0.1739382112039596
$ws.Range("J8").Value = 0.1739382112039596
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 179.590866
$ws.Range("N8").Value = 538.772598
$ws.Range("O8").Value = 0.8800552629635879
$ws.Range("P8").Value = 0.8800552629635878
$ws.Range("Q8").Value = 5368.75729341497
$ws.Range("R8").Value = 48318.81564073473
$ws.Range("S8").Value = 0.1530752382005167
$ws.Range("T8").Value = 0.1530752382005167

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.89437833333333
$ws.Range("H9").Value = 89.68313499999999
$ws.Range("I9").Value = 0.1739382112039596
$ws.Range("J9").Value = 0.1739382112039596
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.589601
$ws.Range("N9").Value = 4.768803
$ws.Range("O9").Value = 0.007789576147275677
$ws.Range("P9").Value = 0.007789576147275675
$ws.Range("Q9").Value = 47.520133693045
$ws.Range("R9").Value = 427.681203237405
$ws.Range("S9").Value = 0.001354904941094163
$ws.Range("T9").Value = 0.001354904941094162

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.711532999999999
$ws.Range("H10").Value = 17.134599
$ws.Range("I10").Value = 0.03323212886968274
$ws.Range("J10").Value = 0.03323212886968274
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.03403033333333
$ws.Range("N10").Value = 63.102091
$ws.Range("O10").Value = 0.1030737782409588
$ws.Range("P10").Value = 0.1030737782409588
$ws.Range("Q10").Value = 120.1365583718343
$ws.Range("R10").Value = 1081.229025346509
$ws.Range("S10").Value = 0.003425361081588644
$ws.Range("T10").Value = 0.003425361081588644

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 5.711532999999999
$ws.Range("H11").Value = 17.134599
$ws.Range("I11").Value = 0.03323212886968274
$ws.Range("J11").Value = 0.03323212886968274
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.853217
$ws.Range("N11").Value = 5.559651000000001
$ws.Range("O11").Value = 0.009081382648177616
$ws.Range("P11").Value = 0.009081382648177616
$ws.Range("Q11").Value = 10.584710051661
$ws.Range("R11").Value = 95.26239046494899
$ws.Range("S11").Value = 0.0003017936784791392
$ws.Range("T11").Value = 0.0003017936784791393

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 5.711532999999999
$ws.Range("H12").Value = 17.134599
$ws.Range("I12").Value = 0.03323212886968274
$ws.Range("J12").Value = 0.03323212886968274
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 179.590866
$ws.Range("N12").Value = 538.772598
$ws.Range("O12").Value = 0.8800552629635879
$ws.Range("P12").Value = 0.8800552629635878
$ws.Range("Q12").Value = 1025.739157657578
$ws.Range("R12").Value = 9231.652418918202
$ws.Range("S12").Value = 0.02924610991124848
$ws.Range("T12").Value = 0.02924610991124849

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 5.711532999999999
$ws.Range("H13").Value = 17.134599
$ws.Range("I13").Value = 0.03323212886968274
$ws.Range("J13").Value = 0.03323212886968274
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.589601
$ws.Range("N13").Value = 4.768803
$ws.Range("O13").Value = 0.007789576147275677
$ws.Range("P13").Value = 0.007789576147275675
$ws.Range("Q13").Value = 9.079058568332998
$ws.Range("R13").Value = 81.711527114997
$ws.Range("S13").Value = 0.000258864198366472
$ws.Range("T13").Value = 0.000258864198366472

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.127415333333333
$ws.Range("H14").Value = 27.382246
$ws.Range("I14").Value = 0.05310718551472111
$ws.Range("J14").Value = 0.05310718551472112
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 21.03403033333333
$ws.Range("N14").Value = 63.102091
$ws.Range("O14").Value = 0.1030737782409588
$ws.Range("P14").Value = 0.1030737782409588
$ws.Range("Q14").Value = 191.9863309862651
$ws.Range("R14").Value = 1727.876978876386
$ws.Range("S14").Value = 0.005473958262745824
$ws.Range("T14").Value = 0.005473958262745825

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.127415333333333
$ws.Range("H15").Value = 27.382246
$ws.Range("I15").Value = 0.05310718551472111
$ws.Range("J15").Value = 0.05310718551472112
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.853217
$ws.Range("N15").Value = 5.559651000000001
$ws.Range("O15").Value = 0.009081382648177616
$ws.Range("P15").Value = 0.009081382648177616
$ws.Range("Q15").Value = 16.915081261794
$ws.Range("R15").Value = 152.235731356146
$ws.Range("S15").Value = 0.0004822866730269379
$ws.Range("T15").Value = 0.0004822866730269381

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.127415333333333
$ws.Range("H16").Value = 27.382246
$ws.Range("I16").Value = 0.05310718551472111
$ws.Range("J16").Value = 0.05310718551472112
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 179.590866
$ws.Range("N16").Value = 538.772598
$ws.Range("O16").Value = 0.8800552629635879
$ws.Range("P16").Value = 0.8800552629635878
$ws.Range("Q16").Value = 1639.200424055012
$ws.Range("R16").Value = 14752.80381649511
$ws.Range("S16").Value = 0.04673725811341393
$ws.Range("T16").Value = 0.04673725811341394

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.127415333333333
$ws.Range("H17").Value = 27.382246
$ws.Range("I17").Value = 0.05310718551472111
$ws.Range("J17").Value = 0.05310718551472112
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.589601
$ws.Range("N17").Value = 4.768803
$ws.Range("O17").Value = 0.007789576147275677
$ws.Range("P17").Value = 0.007789576147275675
$ws.Range("Q17").Value = 14.508948541282
$ws.Range("R17").Value = 130.580536871538
$ws.Range("S17").Value = 0.0004136824655344159
$ws.Range("T17").Value = 0.0004136824655344159
